$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43
$srcRow = 42

# Replicate the formatting of the previous data row (borders/bold index
# column style + date-time number format) onto the new row before writing
# values, so the new cells pick up the same cellXfs entries Excel would
# reuse for a row that was extended in-place.
$ws.Range("A" + $srcRow + ":V" + $srcRow).Copy()
$ws.Range("A" + $row + ":V" + $row).PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 42
$ws.Cells.Item($row, 2).Value = "india"
$ws.Cells.Item($row, 3).Value = "isl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45262.64583333334
$ws.Cells.Item($row, 6).Value = "Hyderabad"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Mohun Bagan"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 2.56
$ws.Cells.Item($row, 11).Value = "25/11/2023 18:13"
$ws.Cells.Item($row, 12).Value = 4.4
$ws.Cells.Item($row, 13).Value = "02/12/2023 15:21"
$ws.Cells.Item($row, 14).Value = 3.09
$ws.Cells.Item($row, 15).Value = "25/11/2023 18:13"
$ws.Cells.Item($row, 16).Value = 3.6
$ws.Cells.Item($row, 17).Value = "02/12/2023 15:21"
$ws.Cells.Item($row, 18).Value = 2.95
$ws.Cells.Item($row, 19).Value = "25/11/2023 18:13"
$ws.Cells.Item($row, 20).Value = 1.85
$ws.Cells.Item($row, 21).Value = "02/12/2023 15:21"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/india/isl/hyderabad-mohun-bagan/CODmlmFL/"
